$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").EntireColumn.Insert()
$ws.Range("A1").Value = "ssim_dual"
Write-Host $ws.Range("A1").Value()
